$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$conv = $wb.Worksheets.Item("CONVERTION")

# --- Insert a new leave-card row above row 91 (shifts 91..130 -> 92..131) ---
$ws.Rows.Item(91).Insert()

# Row 91 lost its table formatting/formula on insert; restore it from the row
# that is now directly below (row 92, the old row 91) so the style ids match
# what a normal in-table row insert produces.
$ws.Range("A92:K92").Copy()
$ws.Range("A91:K91").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(91,7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Column K on row 91 holds a plain date (style like K89), not leave text.
$ws.Cells.Item(89,11).Copy()
$ws.Cells.Item(91,11).PasteSpecial(-4122) | Out-Null

# --- Grow Table1 to include the newly-shifted last row (131) ---
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K131"))
# Re-apply the calculated-column formula on the new last row so the cached
# result isn't left stale/erroring from the resize.
$ws.Cells.Item(131,7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- Fill in the new leave entries ---
$ws.Cells.Item(90,2).Value = "VL(4-0-0)"
$ws.Cells.Item(90,4).Value = 4
$ws.Cells.Item(90,11).Value = "5/16,17,18,19/2023"

$ws.Cells.Item(91,2).Value = "SP(1-0-00)"
$ws.Cells.Item(91,11).Value = 45066

$ws.Cells.Item(92,2).Value = "VL(5-0-00)"
$ws.Cells.Item(92,4).Value = 5
$ws.Cells.Item(92,11).Value = "5/8-12/2023"

# --- View/selection bookkeeping to mirror where the editor ended up ---
$ws.Activate()
$ws.Range("A87").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B93").Select()

$conv.Activate()
$conv.Range("A8").Select()

$ws.Activate()
$ws.Range("B2:C2").Select()
